$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update G8 text
$ws.Range("G8").Value = "无锡一疗住院（4.12）期间"

# Add new row 13
$ws.Range("A13").Value = 45431
$ws.Range("A13").NumberFormat = "m/d/yyyy"
$ws.Range("B13").Value = 1.29
$ws.Range("C13").Value = 0.59
$ws.Range("D13").Value = 1.65
$ws.Range("E13").Value = 60
$ws.Range("F13").Value = 58
$ws.Range("G13").Value = "无锡二疗住院（5.18）期间"
